# Weekly fruit/vegetable price update: add 3 new price records for
# "Macroferia Regional de Talca" - Tomate (Hortaliza), inserted as new
# rows 386-388, pushing the previously-existing rows 386-402 down to
# rows 389-405.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 386, shifting everything
# from 386 downward (old 386-402) down to 389-405.
$ws.Range("A386:R388").EntireRow.Insert()

# --- New row 386 ---------------------------------------------------
$ws.Range("A386").Value = 5
$ws.Range("B386").Value = "Macroferia Regional de Talca"
$ws.Range("C386").Value = "Maule"
$ws.Range("D386").Value = 44509
$ws.Range("E386").Value = 7
$ws.Range("F386").Value = 100112020
$ws.Range("G386").Value = "Tomate"
$ws.Range("H386").Value = "Larga vida"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 2500
$ws.Range("K386").Value = 13000
$ws.Range("L386").Value = 13000
$ws.Range("M386").Value = 13000
$ws.Range("N386").Value = "`$/bandeja 15 kilos"
$ws.Range("O386").Value = "Región del Maule"
$ws.Range("P386").Value = 867
$ws.Range("Q386").Value = 15
$ws.Range("R386").Value = "Hortaliza"

# --- New row 387 ---------------------------------------------------
$ws.Range("A387").Value = 5
$ws.Range("B387").Value = "Macroferia Regional de Talca"
$ws.Range("C387").Value = "Maule"
$ws.Range("D387").Value = 44509
$ws.Range("E387").Value = 7
$ws.Range("F387").Value = 100112020
$ws.Range("G387").Value = "Tomate"
$ws.Range("H387").Value = "Larga vida"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 2500
$ws.Range("K387").Value = 15000
$ws.Range("L387").Value = 15000
$ws.Range("M387").Value = 15000
$ws.Range("N387").Value = "`$/bandeja 18 kilos"
$ws.Range("O387").Value = "Región de Arica y Parinacota"
$ws.Range("P387").Value = 833
$ws.Range("Q387").Value = 18
$ws.Range("R387").Value = "Hortaliza"

# --- New row 388 ---------------------------------------------------
$ws.Range("A388").Value = 5
$ws.Range("B388").Value = "Macroferia Regional de Talca"
$ws.Range("C388").Value = "Maule"
$ws.Range("D388").Value = 44509
$ws.Range("E388").Value = 7
$ws.Range("F388").Value = 100112020
$ws.Range("G388").Value = "Tomate"
$ws.Range("H388").Value = "Larga vida"
$ws.Range("I388").Value = "Primera"
$ws.Range("J388").Value = 2000
$ws.Range("K388").Value = 18000
$ws.Range("L388").Value = 18000
$ws.Range("M388").Value = 18000
$ws.Range("N388").Value = "`$/bandeja 18 kilos"
$ws.Range("O388").Value = "Región del Maule"
$ws.Range("P388").Value = 1000
$ws.Range("Q388").Value = 18
$ws.Range("R388").Value = "Hortaliza"
